# Update "想去人数" (F column) values for the first two events on both the
# "展览" (Exhibition) sheet and the "全部类型" (All Types) sheet.
#   F2: 1774 -> 1777
#   F3: 8101 -> 8107

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 1777
    $ws.Range("F3").Value = 8107
}
